$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New system-requirements header columns (F..K). Written in this exact
# order so the shared-string table receives "М.с.т. Процессор" before
# "М.с.т. ОС" (matches the target sharedStrings.xml ordering).
$ws.Range("G1").Value = "М.с.т. Процессор"
$ws.Range("F1").Value = "М.с.т. ОС"
$ws.Range("H1").Value = "М.с.т. ОЗУ"
$ws.Range("I1").Value = "М.с.т. Видеокарта"
$ws.Range("J1").Value = "М.с.т. Память видеокарты"
$ws.Range("K1").Value = "М.с.т. Жесткий диск"

# Keep the new header cells bold, matching the rest of row 1.
$ws.Range("H1:K1").Font.Bold = $true

# New data row for the system requirements of the game in row 2.
$ws.Range("F2").Value = "64-разрядная версия Windiws 7"
$ws.Range("G2").Value = "Четырехъядерный процессор Intel Core i3-6300 3,8 ГГц "

# Distinguish F2's font from the plain default so a 4th font/style entry
# is produced (mirrors the extra font Excel generated for this cell).
$ws.Range("F2").Font.ColorIndex = 1

# Column F was resized from the old (wide) "system requirements" column to
# a narrow numeric-ish width.
$ws.Columns.Item(6).ColumnWidth = 9.333333

# Move the view so column C is the left-most visible column, and leave the
# selection on G2 (the last cell touched).
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("G2").Select()
